# Refresh the cryptocurrency ranking snapshot (prices / 1h volume%, and a few
# rank swaps between neighbouring coins) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as literal text, never letting Excel
# auto-convert number-shaped strings (e.g. "1.00", "75.237.79", "0.0000192")
# into real numbers / scientific notation.
function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# Row 2
Set-TextCell 2 4 "75.290.18"
Set-TextCell 2 5 "  +1.82%  "
# Row 3
Set-TextCell 3 4 "2.869.38"
Set-TextCell 3 5 "  +10.99%  "
# Row 4
Set-TextCell 4 5 "  +0.09%  "
# Row 5
Set-TextCell 5 4 "608.90"
Set-TextCell 5 5 "  +4.54%  "
# Row 6
Set-TextCell 6 4 "189.90"
Set-TextCell 6 5 "  +4.87%  "
# Row 7
Set-TextCell 7 5 "  -0.06%  "
# Row 8
Set-TextCell 8 4 "0.572"
Set-TextCell 8 5 "  +7.55%  "
# Row 9
Set-TextCell 9 4 "0.196"
Set-TextCell 9 5 "  -4.73%  "
# Row 10
Set-TextCell 10 4 "2.873.34"
Set-TextCell 10 5 "  +11.05%  "
# Row 11
$ws.Cells.Item(11, 2).Value = "Cardano"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextCell 11 4 "0.378"
Set-TextCell 11 5 "  +6.60%  "
# Row 12
$ws.Cells.Item(12, 2).Value = "TRON"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell 12 4 "0.163"
Set-TextCell 12 5 "  +0.62%  "
# Row 13
Set-TextCell 13 4 "4.99"
Set-TextCell 13 5 "  +4.35%  "
# Row 14
Set-TextCell 14 4 "3.394.71"
Set-TextCell 14 5 "  +11.30%  "
# Row 15
$ws.Cells.Item(15, 2).Value = "Avalanche"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell 15 4 "27.94"
Set-TextCell 15 5 "  +7.43%  "
# Row 16
$ws.Cells.Item(16, 2).Value = "WrappedBTC"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextCell 16 4 "75.237.79"
Set-TextCell 16 5 "  +2.13%  "
# Row 17
$ws.Cells.Item(17, 2).Value = "ShibaInu"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell 17 4 "0.0000192"
Set-TextCell 17 5 "  -0.31%  "
# Row 18
Set-TextCell 18 4 "2.874.74"
Set-TextCell 18 5 "  +12.06%  "
# Row 19
Set-TextCell 19 4 "9.31"
Set-TextCell 19 5 "  +17.77%  "
# Row 20
Set-TextCell 20 5 "  +8.71%  "
# Row 21
Set-TextCell 21 4 "382.65"
Set-TextCell 21 5 "  +5.05%  "
# Row 22
Set-TextCell 22 5 "  +3.58%  "
# Row 23
Set-TextCell 23 4 "4.19"
Set-TextCell 23 5 "  +3.41%  "
# Row 24
Set-TextCell 24 4 "6.24"
Set-TextCell 24 5 "  +0.43%  "
# Row 25
Set-TextCell 25 4 "71.70"
Set-TextCell 25 5 "  +3.90%  "
# Row 26
Set-TextCell 26 5 "  +0.07%  "
# Row 27
Set-TextCell 27 4 "4.31"
Set-TextCell 27 5 "  +4.54%  "
# Row 28
Set-TextCell 28 4 "3.003.14"
Set-TextCell 28 5 "  +10.84%  "
# Row 29
Set-TextCell 29 4 "9.76"
Set-TextCell 29 5 "  +7.68%  "
# Row 30
Set-TextCell 30 4 "0.0000106"
Set-TextCell 30 5 "  +13.74%  "
# Row 31
Set-TextCell 31 4 "0.998"
Set-TextCell 31 5 "  -0.37%  "
# Row 32
Set-TextCell 32 4 "539.45"
Set-TextCell 32 5 "  +8.67%  "
# Row 33
Set-TextCell 33 4 "1.43"
Set-TextCell 33 5 "  +7.70%  "
# Row 34
Set-TextCell 34 4 "8.01"
Set-TextCell 34 5 "  +1.83%  "
# Row 35
Set-TextCell 35 5 "  +9.07%  "
# Row 36
Set-TextCell 36 4 "1.00"
Set-TextCell 36 5 "  +0.07%  "
# Row 37
Set-TextCell 37 4 "0.122"
Set-TextCell 37 5 "  +3.41%  "
# Row 38
Set-TextCell 38 4 "20.48"
Set-TextCell 38 5 "  +7.40%  "
# Row 39
Set-TextCell 39 4 "162.48"
Set-TextCell 39 5 "  +1.44%  "
# Row 40
$ws.Cells.Item(40, 2).Value = "WhiteBITCoin"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextCell 40 4 "19.32"
Set-TextCell 40 5 "  -0.04%  "
# Row 41
$ws.Cells.Item(41, 2).Value = "Aave"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell 41 4 "185.69"
Set-TextCell 41 5 "  +25.90%  "
# Row 43
Set-TextCell 43 4 "5.21"
Set-TextCell 43 5 "  +8.23%  "
# Row 44
Set-TextCell 44 4 "0.347"
Set-TextCell 44 5 "  +9.44%  "
# Row 45
Set-TextCell 45 4 "1.72"
Set-TextCell 45 5 "  +3.47%  "
# Row 46
Set-TextCell 46 5 "  +11.91%  "
# Row 47
Set-TextCell 47 4 "2.41"
Set-TextCell 47 5 "  +1.20%  "
# Row 48
Set-TextCell 48 4 "40.02"
Set-TextCell 48 5 "  +1.94%  "
# Row 49
Set-TextCell 49 4 "0.0861"
Set-TextCell 49 5 "  +8.39%  "
# Row 50
Set-TextCell 50 4 "0.586"
Set-TextCell 50 5 "  +13.25%  "
# Row 51
Set-TextCell 51 4 "3.82"
Set-TextCell 51 5 "  +7.32%  "
